# Remove the trailing "blank line / page-break / copyright footer" block
# that was appended to the bibliography section, restoring the document to
# end right after the last bibliography entry (with its original closing
# blank paragraph + page-break paragraph kept intact).

$d = $word.ActiveDocument

# Locate the paragraph that holds the copyright/footer text and use it as
# an anchor, so the removal is robust even if paragraph indices shift.
$footerText = "Contact: luizeleno@usp.br"

$found = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*$footerText*") {
        $footerIndex = $i
        $found = $true
        break
    }
}

if ($found) {
    # The block to delete is the 3 paragraphs ending with the footer
    # paragraph itself:
    #   1) an empty "Normal" paragraph
    #   2) an empty paragraph with a page break before it
    #   3) the paragraph containing the copyright/footer text
    $startIndex = $footerIndex - 2
    $startPara = $d.Paragraphs($startIndex)
    $endPara = $d.Paragraphs($footerIndex)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
